$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "hussein1992024082300"
$ws.Range("B3").Value = "Test@123"
$ws.Range("A4").Value = "hussein1992024082712"
$ws.Range("B4").Value = "Test@123"
